$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure the tree: add a new "D" row/column to the pairwise comparison matrix ---

# New header in D1
$ws.Range("D1").Value = "D"

# New value in D2 (row "B")
$ws.Range("D2").Value = 0.2

# Row 3 ("C") gets its C3 value kept at 1 and a new D3 value
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 0.2

# New row 4 ("D")
$ws.Range("A4").Value = "D"
$ws.Range("B4").Value = 5
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 1

# C2 previously carried a stray explicit number-format style; clear it back to default
$ws.Range("C2").ClearFormats()

# Match the author's last selection in the sheet
[void]$ws.Range("C5").Select()
